# Update the shadow-price/model results in column D of Sheet1
# to reflect the corrected ("fixed code") model run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value  = 3053.371704751997
$ws.Range("D3").Value  = 3053.371704751994
$ws.Range("D4").Value  = 3053.371704751994

$ws.Range("D6").Value  = 8723.919156434198
$ws.Range("D7").Value  = 8723.919156434198

$ws.Range("D9").Value  = 6876.196583938373
$ws.Range("D10").Value = 6876.196583938367
$ws.Range("D11").Value = 100
$ws.Range("D12").Value = 5676.196583938367
$ws.Range("D13").Value = 100

$ws.Range("D17").Value = 80000

$ws.Range("D19").Value = 89076.78952001187
$ws.Range("D20").Value = 89076.78952001187
$ws.Range("D21").Value = 1200

$ws.Range("D24").Value = 148176.208985743
$ws.Range("D25").Value = 148176.2089857429

$ws.Range("D28").Value = 2963.524179714889
$ws.Range("D29").Value = 2963.524179714889
$ws.Range("D30").Value = 148176.2089857429

$ws.Range("D38").Value = -3871.465396434346
$ws.Range("D39").Value = -3871.465396434342

$ws.Range("D41").Value = 3871.465396434342
$ws.Range("D42").Value = 77429.30792868507
$ws.Range("D43").Value = 77429.30792868507
